$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S1").Value = "DD"
$ws.Range("T1").Value = "DD"

$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 2

$ws.Range("S3").Value = 3
$ws.Range("T3").Value = 4

$ws.Range("S2").Select()
